# Build site at 2023-04-12 14:53:07 UTC
# Updates the LOT2017 "Enzimologia" syllabus sheet:
#  - Rewrites the "Objetivos:" (objectives) paragraph in row 10
#  - Inserts a new row 13 carrying the "Docentes responsaveis:" name
#    (previously mis-placed in row 10), shifting the rest of the sheet down
#  - Rewrites "Programa resumido:", "Programa:" and adds the Bibliografia
#    text block, plus reorders the evaluation paragraphs

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$objectivesPt = "Formação dos estudantes de Engenharia Bioquímica na área de tecnologia de enzimas, com foco principal nos estudos de estrutura versus propriedades e mecanismos de ação, controle operacional na purificação e imobilização de enzimas, formas de determinação de atividade enzimática e aplicações das enzimas nos processos industriais."

$professorLine = "4873328 - Fernando Segato"

$shortProgramPt = "A disciplina aborda como as enzimas atuam, como se definem as estratégias de purificação e quais são as principais aplicações tecnológicas das enzimas. Dentro dos processos de purificação, o foco envolve a definição de estratégias apropriadas para a purificação em etapas sequenciais, os métodos de controle de cada etapa, além dos métodos de monitoramento da atividade enzimática. Também se aborda a aplicação das enzimas em processos industriais."

$fullProgramPt = "1. Origem celular das enzimas: origem das enzimas, diferenciação entre enzimas intra e extracelulares, importância fisiológica e introdução ao mercado mundial de enzimas.2. Estrutura versus propriedades e mecanismos de ação das enzimas: estruturas tridimensionais e sua determinação, importância da estrutura terciária na atividade catalítica, ação catalítica de proteases, glicosidases e oxido-redutases.3. Controle operacional na purificação de enzimas: métodos de extração de enzimas, métodos de purificação preliminar, métodos de separação baseados na carga, no tamanho e na afinidade. Definição de estratégias de purificação de enzimas.4. Métodos de determinação de atividade enzimática: definição de atividade em unidades internacionais, significado da atividade enzimática, formas de quantificar e expressar e atividade enzimática. Requerimentos de um método experimental usado na determinação de atividade enzimática.5. Cinética enzimática: métodos gráficos e numéricos de determinação de velocidade inicial de reação, condições experimentais demandadas para determinar a velocidade inicial, cálculos de atividade enzimática.6. Enzimas imobilizadas: formas de imobilização e aplicações de sistemas imobilizados.7. Aplicações de enzimas na indústria: uso de enzimas em detergentes, no processamento do amido, na indústria alimentícia, na indústria têxtil, na síntese de fármacos e na indústria de celulose e papel."

$criterioPt = "A avaliação será feita por meio de provas escritas (P1 e P2)."
$normaRecPt = "A Nota final (NF) será calculada da seguinte maneira: NF = (P1x1 + P2x2)/3"
$bibliografiaIntro = "A recuperação será feita por meio de uma prova escrita (PR) e a média de recuperação (MR) calculada pela fórmula: MR = (NF + PR)/2"

$bibliografiaText = "1. BON, E.S., FERRARA M.A., CORVO M.L. (Eds.) Enzimas em Biotecnologia - Produção, aplicação e mercado, Rio de Janeiro: Editora Interciêcnia, 2008.`n2. COPELAND, R.A. Enzymes: a practical introduction to structure, mechanism and data analysis, New York: Academic Press, 2000.`n3. LEHNINGER, A.L., NELSON, O.L., COX, M.M. Princípios de bioquímica, 5 ed. Porto Alegre: Artmed editora, 2011.`n4. GODFREY, T., WEST, S. (eds), Industrial Enzymology, New York: Chapman-Hall, 1996.`n5. WHITAKER, J.R. (ed.) Pinciples of Enzynmology for the Food Sciences 2nd ed., New York: Marcel Dekker Inc., 1994.`n6. TANAKA, A., TOSA, T., KOBAYASHI, T. (Eds.). Industrial Application of Immobilized Biocatalysts, New York: Marcel Dekker Inc., 1993.`n6.VOET, D., VOET, J., PRATT, C.W. Fundamentos de Bioquímica. Porto Alegre: Editora ARTMED, 2000."

# 1) Row 10 ("Objetivos:") - replace the professor name with the real objectives text
$ws.Range("B10").Value = $objectivesPt
$ws.Range("C10").Value = $objectivesPt

# 2) Insert a new row at 13 (pushes everything from the old row 13 onward down by one)
$ws.Rows(13).Insert()

# New row 13 holds the professor line under "Docentes responsaveis:" (row 12).
# It needs no value in column A, so copy the B/C number formats from row 10
# (same column styles) before writing the values.
$ws.Range("B10:C10").Copy()
$ws.Range("B13").PasteSpecial(-4122)
$ws.Range("A13").Clear()
$ws.Range("B13").Value = $professorLine
$ws.Range("C13").Value = $professorLine

# 3) Row 14 ("Programa resumido:", was old row 13) - replace "Semestral" with the real summary
$ws.Range("B14").Value = $shortProgramPt
$ws.Range("C14").Value = $shortProgramPt

# 4) Row 16 ("Programa:", was old row 15) - replace the placeholder date with the full program text
$ws.Range("B16").Value = $fullProgramPt
$ws.Range("C16").Value = $fullProgramPt

# 5) Row 19 ("Metodo:", was old row 18) - replace the professor name with the evaluation method text
$ws.Range("B19").Value = $criterioPt
$ws.Range("C19").Value = $criterioPt

# 6) Row 20 ("Criterio:", was old row 19) - shift in the final-grade formula text
$ws.Range("B20").Value = $normaRecPt
$ws.Range("C20").Value = $normaRecPt

# 7) Row 21 ("Norma de recuperacao:", was old row 20) - shift in the recovery exam text
$ws.Range("B21").Value = $bibliografiaIntro
$ws.Range("C21").Value = $bibliografiaIntro

# 8) Row 22 ("Bibliografia:", was old row 21) - add the bibliography list (new cells, copy format from row 21)
$ws.Range("B21:C21").Copy()
$ws.Range("B22:C22").PasteSpecial(-4122)
$ws.Range("B22").Value = $bibliografiaText
$ws.Range("C22").Value = $bibliografiaText

Write-Host "Edit complete"
